# Adds "salesTest" family + "priceUpdateTest" rows to both sheets, right
# before the trailing "chainSummaryTest" row, mirroring the existing
# pattern of test-case rows. Matches commit: "price update, reservation
# and sales test is added".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Sheet1"
$ws2 = $wb.Worksheets.Item(2)   # "Data"

$names = "salesTest","salesTestForLastYearOff","salesTestForLastYearOffDownload","salesTestDownload","priceUpdateTest"
$execs = "no","no","no","no","yes"

# ---------------------------------------------------------------------
# Sheet1 ("Sheet1"): columns A..E = TestCaseName / Description / Execute /
# InvocationCount / Priority
# ---------------------------------------------------------------------

# Insert 5 fresh rows right before the last row (old row 13), pushing it
# down to row 18, then stamp the formatting of row 12 onto them so the
# look (borders/fill/font) matches the rest of the table exactly.
$ws1.Range("A13:E17").Insert(-4121)
$ws1.Range("A12:E12").Copy()
$ws1.Range("A13:E17").PasteSpecial(-4122)

for ($i = 0; $i -lt 5; $i++) {
    $r = 13 + $i
    $ws1.Range("A$r").Value = $names[$i]
    $ws1.Range("B$r").Value = "abcd"
    $ws1.Range("C$r").Value = $execs[$i]
    $ws1.Range("D$r").Value = "'1"
    $ws1.Range("E$r").Value = "'1"
}

# Former last test-case row (row 12, "optimizerTest") is no longer the
# final row, so its Execute flag flips from "yes" to "no".
$ws1.Range("C12").Value = "no"

# The "no data past here" duplicate-check range shifts down with the
# newly inserted rows.
$fcs1 = $ws1.Range("A1:A1048576").FormatConditions
$fcs1.Item(2).ModifyAppliesToRange($ws1.Range("A19:A1048576"))

# ---------------------------------------------------------------------
# Sheet2 ("Data"): columns A..E = TestName / Execute / Browser / UserName
# / Password
# ---------------------------------------------------------------------

$ws2.Range("A13:E17").Insert(-4121)
$ws2.Range("A12:E12").Copy()
$ws2.Range("A13:E17").PasteSpecial(-4122)

for ($i = 0; $i -lt 5; $i++) {
    $r = 13 + $i
    $ws2.Range("A$r").Value = $names[$i]
    $ws2.Range("B$r").Value = $execs[$i]
    $ws2.Range("C$r").Value = "chrome"
    $ws2.Range("D$r").Value = "raghavendra.m@axisrooms.com"
    $ws2.Range("E$r").Value = "Password123#"
}

$ws2.Range("B12").Value = "no"

$fcs2 = $ws2.Range("A1:A1048576").FormatConditions
$fcs2.Item(2).ModifyAppliesToRange($ws2.Range("A19:A1048576"))
$fcs2.Item(3).ModifyAppliesToRange($ws2.Range("A19:A1048576"))

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping: selection on both sheets moves to
# the new last data row (A17), and the "Data" tab ends up the active one.
# ---------------------------------------------------------------------

$ws1.Activate()
$ws1.Range("A17").Select()

$ws2.Activate()
$ws2.Range("A17").Select()
